$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.539.38"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.824.16"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'315.64"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5106"
$ws.Range("E7").Value = "  -5.51%  "
$ws.Range("D8").Value = "'0.3953"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").Value = "'0.08246"
$ws.Range("E9").Value = "  +6.19%  "
$ws.Range("D10").Value = "'1.113"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").Value = "'41.71"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "'6.365"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "'21.18"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "'1.002"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "'7.555"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "1.819.25"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "'0.00001126"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "'93.03"
$ws.Range("E18").Value = "  +3.38%  "
$ws.Range("D19").Value = "'0.06656"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").Value = "'17.84"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'6.106"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "28.585.71"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'11.43"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").Value = "'2.276"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").Value = "'21.34"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").Value = "2.030.45"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "'2.418"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").Value = "'126.68"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "'5.785"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "'3.662"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").Value = "'0.07062"
$ws.Range("E35").Value = "  -6.20%  "
$ws.Range("D36").Value = "'0.2231"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").Value = "'0.02355"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "'5.265"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Value = "'8.799"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "'0.6354"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").Value = "'1.405"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "'13.59"
$ws.Range("D45").Value = "'0.5964"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").Value = "'3.737"
$ws.Range("D47").Value = "'125.29"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "'1.196"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'0.06944"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "'1.082"
